$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) and "Volume(1h)" (E) columns hold values that look numeric
# (e.g. "42.979.75", "1.00", "34.85") but must stay plain text, exactly as
# authored upstream (inline strings). Excel's Range.Value setter auto-detects
# number-like text and silently coerces it to a real number, which would
# corrupt values like "1.00" -> 1. To avoid that, temporarily force the
# target range to Text format, write the values, then restore the original
# style so no stray formatting diff is introduced.
$dataRange = $ws.Range("B2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.003.80"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.339.34"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "307.43"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").Value = "101.44"
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("D7").Value = "0.510"
$ws.Range("E7").Value = "  -4.35%  "
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  -3.53%  "
$ws.Range("D10").Value = "34.85"
$ws.Range("E10").Value = "  -4.47%  "
$ws.Range("D11").Value = "52.45"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").Value = "0.0798"
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").Value = "6.86"
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").Value = "15.83"
$ws.Range("E15").Value = "  +5.08%  "
$ws.Range("D16").Value = "2.348.47"
$ws.Range("E16").Value = "  +5.63%  "
$ws.Range("D17").Value = "0.832"
$ws.Range("E17").Value = "  +2.59%  "
$ws.Range("D18").Value = "42.919.43"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").Value = "  -4.92%  "
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "69.04"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "236.74"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("D24").Value = "2.02"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "25.60"
$ws.Range("E27").Value = "  +2.90%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "3.96"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.32"
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "35.29"
$ws.Range("E30").Value = "  -4.32%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Value = "9.30"
$ws.Range("E31").Value = "  -4.10%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "161.86"
$ws.Range("E32").Value = "  -3.79%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "5.12"
$ws.Range("E34").Value = "  -3.30%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "4.67"
$ws.Range("E35").Value = "  +5.29%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "17.41"
$ws.Range("E36").Value = "  -3.66%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "2.45"
$ws.Range("E37").Value = "  -3.67%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.0726"
$ws.Range("E38").Value = "  -2.40%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "1.85"
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "2.91"
$ws.Range("E40").Value = "  -5.18%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.102"
$ws.Range("E41").Value = "  -3.75%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.113"
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "2.61"
$ws.Range("E43").Value = "  +4.41%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.021.29"
$ws.Range("E44").Value = "  +1.82%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0285"
$ws.Range("E45").Value = "  -3.72%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "18.98"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "10.24"
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "2.94"
$ws.Range("E48").Value = "  -3.02%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "55.94"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "2.89"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.563.88"
$ws.Range("E51").Value = "  +1.04%  "

$dataRange.Style = $origStyle
